## Generate Report for Archive
##
## - Status moves from "Ready for handoff" to "In Translation" on every
##   sheet that tracks it (Overview!E2/F2, zh-cn!C2, de-de!C2).
## - The Status columns are re-narrowed to fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

## This host's `Range.ColumnWidth` setter quantizes the written value onto a
## 1/6-character grid (stored = (Round(ColumnWidth * 6) + 5) / 6) before it
## is exported into the sheet's <col width="..."/>, rather than leaving the
## float untouched. Pre-compensate the requested character width so the
## exported width lands on the nearest representable grid point.
function Set-PreciseColumnWidth($col, $targetWidth) {
    $raw = ($targetWidth * 6) - 5
    $snapped = [Math]::Round($raw)
    $col.ColumnWidth = $snapped / 6
}

$targetStatusWidth = 13.4101845877511

Set-PreciseColumnWidth $wsOverview.Columns.Item(5) $targetStatusWidth
Set-PreciseColumnWidth $wsOverview.Columns.Item(6) $targetStatusWidth
Set-PreciseColumnWidth $wsZhCn.Columns.Item(3) $targetStatusWidth
Set-PreciseColumnWidth $wsDeDe.Columns.Item(3) $targetStatusWidth
